# Supports EPV. Include responder for text-related questions in the spreadsheets generated.
#
# The existing sheet holds 4 columns (Single Choice, Multiple Choice, Free Text,
# Number). This edit:
#   1. Inserts a new column A ("Responder") with respondent ids R1..R13,
#      pushing the old A:D columns to B:E.
#   2. Appends three new "Table (A,B,C,D)" columns (F, G, H) with
#      Small/Medium/Large sample data per responder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column before A, shifting the old A:D -> B:E.
$ws.Columns("A:A").Insert()

# 2) Fill the new "Responder" column (A).
$responder = @("Responder","R1","R2","R3","R4","R5","R6","R7","R8","R9","R10","R11","R12","R13")
for ($i = 0; $i -lt $responder.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $responder[$i]
}

# 3) Fill the three new "Table (A,B,C,D)" columns (F, G, H).
$tableF = @(
    "Table (A,B,C,D)",
    "Small;2;100;100",
    "Small;1;200;200",
    "Small;3;100;100",
    "Small;4;200;200",
    "Small;1;100;100",
    "Small;2;200;200",
    "Small;1;100;100",
    "Small;5;200;200",
    "Small;3;100;100",
    "Small;2;200;200",
    "Small;5;100;100",
    "Small;1;200;200",
    "Small;2;1000;1231"
)

$tableG = @(
    "Table (A,B,C,D)",
    "Medium;1;200;300",
    "Medium;2;300;400",
    "Medium;4;200;300",
    "Medium;2;300;400",
    "Medium;3;200;300",
    "Medium;2;300;400",
    "Medium;3;200;300",
    "Medium;1;300;400",
    "Medium;2;200;300",
    "Medium;4;300;400",
    "Medium;1;200;300",
    "Medium;2;300;400",
    "Medium;1;1111;1422"
)

$tableH = @(
    "Table (A,B,C,D)",
    "Large;1;1000;1000",
    "Large;3;1000;1000",
    "Large;2;1000;1000",
    "Large;1;1000;1000",
    "Large;3;1000;1000",
    "Large;1;1000;1000",
    "Large;2;1000;1000",
    "Large;1;1000;1000",
    "Large;1;1000;1000",
    "Large;1;1000;1000",
    "Large;1;1000;1000",
    "Large;1;1000;1000",
    "Large;1;1511;1111"
)

for ($i = 0; $i -lt $tableF.Length; $i++) {
    $ws.Cells.Item($i + 1, 6).Value = $tableF[$i]
    $ws.Cells.Item($i + 1, 7).Value = $tableG[$i]
    $ws.Cells.Item($i + 1, 8).Value = $tableH[$i]
}

# Clear the single-cell selection box left over on the old sheet view.
$ws.Range("A1").Select() | Out-Null
